$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 (tato / 1234) - shifts nothing below it, just removes it
$ws.Rows.Item(5).Delete()

# Update A3/B3: was "admin"/"donpedro", now "tato"/1234
$ws.Range("A3").Value = "tato"
$ws.Range("B3").Value = 1234

# Update B4: javier's password now 5289
$ws.Range("B4").Value = 5289

# Update selection to A5 (empty row below data)
$ws.Range("A5").Select()
